$d = $word.ActiveDocument

$pairs = @(
    @("492×7=3444", "477×4=1908"),
    @("191×9=1719", "232×8=1856"),
    @("800×4=3200", "465×3=1395"),
    @("114×5=570",  "751×5=3755"),
    @("166×4=664",  "315×4=1260"),
    @("931×3=2793", "514×8=4112"),
    @("583×2=1166", "482×6=2892"),
    @("531×2=1062", "429×8=3432"),
    @("634×6=3804", "931×7=6517"),
    @("307×5=1535", "782×8=6256"),
    @("604×5=3020", "595×9=5355"),
    @("614×9=5526", "222×3=666"),
    @("321×9=2889", "150×3=450"),
    @("874×5=4370", "909×9=8181"),
    @("825×6=4950", "672×7=4704"),
    @("946×7=6622", "535×5=2675"),
    @("203×4=812",  "348×2=696"),
    @("629×9=5661", "498×8=3984"),
    @("831×8=6648", "488×8=3904"),
    @("103×3=309",  "211×7=1477"),
    @("845×8=6760", "130×4=520"),
    @("914×2=1828", "206×7=1442"),
    @("253×6=1518", "273×2=546"),
    @("684×9=6156", "923×9=8307"),
    @("738×4=2952", "469×2=938")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
